# The deck's slide-master theme ("Integral") is being swapped for the
# plain default "Office Theme" palette (this mirrors the authored commit,
# which exchanged the contents of ppt/theme/theme1.xml and
# ppt/theme/theme2.xml so the slide master now carries the stock Office
# color scheme). Only the 10 colors that actually differ between the two
# themes need updating; dk1/lt1 (black/white) are identical already.
#
# COM's RGBColor.RGB uses the standard Windows COLORREF (little-endian
# BGR) packing, i.e. 0xBBGGRR for a source hex RRGGBB.

$p  = $ppt.ActivePresentation
$cs = $p.SlideMaster.Theme.ThemeColorScheme

$cs.Item(3).RGB  = 0x6A5444   # dk2      455F51 -> 44546A
$cs.Item(4).RGB  = 0xE6E6E7   # lt2      E3DED1 -> E7E6E6
$cs.Item(5).RGB  = 0xD59B5B   # accent1  99CB38 -> 5B9BD5
$cs.Item(6).RGB  = 0x317DED   # accent2  63A537 -> ED7D31
$cs.Item(7).RGB  = 0xA5A5A5   # accent3  E6D024 -> A5A5A5
$cs.Item(8).RGB  = 0xC0FF     # accent4  CC9700 -> FFC000
$cs.Item(9).RGB  = 0xC47244   # accent5  4EB3CF -> 4472C4
$cs.Item(10).RGB = 0x47AD70   # accent6  378DA6 -> 70AD47
$cs.Item(11).RGB = 0xC16305   # hlink    6B9F25 -> 0563C1
$cs.Item(12).RGB = 0x724F95   # folHlink B26B02 -> 954F72
